$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new (empty) column before column M. This shifts the old
# M, N, O columns (delta/IC_Inf/IC_Sup numeric columns for 70_up) one
# column to the right, becoming N, O, P, matching the target layout.
$ws.Columns("M:M").Insert()

# Widen column B to fit the new longer row labels.
$ws.Columns("B:B").ColumnWidth = 33.77

# Append six new data rows (17-22) for the "elderly" / "non elderly"
# breakdown, copying the number formatting from row 16 first so the
# new cells share the existing styles (no new style entries created).
$ws.Range("A16:P16").Copy()
$ws.Range("A17:P22").PasteSpecial(-4122)

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "admissoes_gerais_non_elderly"
$ws.Range("C17").Value = 5039982
$ws.Range("D17").Value = 4935269
$ws.Range("E17").Value = 4890078
$ws.Range("F17").Value = 4965148
$ws.Range("G17").Value = 4922256
$ws.Range("H17").Value = 4826905
$ws.Range("I17").Value = 4910908
$ws.Range("J17").Value = 5077993
$ws.Range("K17").Value = 5143596
$ws.Range("L17").Value = 2.05584067562146
$ws.Range("N17").Value = 0.250598342006092
$ws.Range("O17").Value = 0.239218167657662
$ws.Range("P17").Value = 0.261979808347523

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "admissoes_gerais_uti_non_elderly"
$ws.Range("C18").Value = 153443
$ws.Range("D18").Value = 156757
$ws.Range("E18").Value = 156712
$ws.Range("F18").Value = 161163
$ws.Range("G18").Value = 160838
$ws.Range("H18").Value = 162158
$ws.Range("I18").Value = 166564
$ws.Range("J18").Value = 171725
$ws.Range("K18").Value = 182960
$ws.Range("L18").Value = 19.2364591411795
$ws.Range("N18").Value = 1.88977113447928
$ws.Range("O18").Value = 1.82600366546395
$ws.Range("P18").Value = 1.95357853720464

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "admissoes_gerais_non_uti_non_elderly"
$ws.Range("C19").Value = 4886539
$ws.Range("D19").Value = 4778512
$ws.Range("E19").Value = 4733366
$ws.Range("F19").Value = 4803985
$ws.Range("G19").Value = 4761418
$ws.Range("H19").Value = 4664747
$ws.Range("I19").Value = 4744344
$ws.Range("J19").Value = 4906268
$ws.Range("K19").Value = 4960636
$ws.Range("L19").Value = 1.51634930162227
$ws.Range("N19").Value = 0.195280386577612
$ws.Range("O19").Value = 0.18371454090611
$ws.Range("P19").Value = 0.206847567483948

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "admissoes_gerais_elderly"
$ws.Range("C20").Value = 2144643
$ws.Range("D20").Value = 2125158
$ws.Range("E20").Value = 2186267
$ws.Range("F20").Value = 2223306
$ws.Range("G20").Value = 2260781
$ws.Range("H20").Value = 2261515
$ws.Range("I20").Value = 2348437
$ws.Range("J20").Value = 2416990
$ws.Range("K20").Value = 2525487
$ws.Range("L20").Value = 17.7579205490145
$ws.Range("N20").Value = 2.04300210867354
$ws.Range("O20").Value = 2.02587807324528
$ws.Range("P20").Value = 2.06012901820192

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "admissoes_gerais_uti_elderly"
$ws.Range("C21").Value = 180067
$ws.Range("D21").Value = 189107
$ws.Range("E21").Value = 199676
$ws.Range("F21").Value = 212552
$ws.Range("G21").Value = 218000
$ws.Range("H21").Value = 222873
$ws.Range("I21").Value = 236761
$ws.Range("J21").Value = 249121
$ws.Range("K21").Value = 269248
$ws.Range("L21").Value = 49.526565111875
$ws.Range("N21").Value = 4.84054849133848
$ws.Range("O21").Value = 4.78371011611165
$ws.Range("P21").Value = 4.89741769770207

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "admissoes_gerais_non_uti_elderly"
$ws.Range("C22").Value = 1964576
$ws.Range("D22").Value = 1936051
$ws.Range("E22").Value = 1986591
$ws.Range("F22").Value = 2010754
$ws.Range("G22").Value = 2042781
$ws.Range("H22").Value = 2038642
$ws.Range("I22").Value = 2111676
$ws.Range("J22").Value = 2167869
$ws.Range("K22").Value = 2256239
$ws.Range("L22").Value = 14.846104197547
$ws.Range("N22").Value = 1.74968914783749
$ws.Range("O22").Value = 1.73172972834583
$ws.Range("P22").Value = 1.76765173783211

# Match the author's final selection/active cell.
[void]$ws.Range("E24").Select()
